# Commit: Mon, May 18, 2020 9:05:27 PM
#
# The underlying XML diff shows a single semantic edit inside the deck:
# the data table on slide 6 ("SOURCES OF FINANCE") had a new Table
# Style applied to it (its <a:tableStyleId> GUID changed from the
# default table style to the "{FA53ABA6-793C-4A8D-89ED-51FAC7A60CF6}"
# gallery style). Reproduce that with the Table Style gallery API
# (PowerPoint COM requires Table.ApplyStyle("{GUID}") -- Table.Style
# itself is read-only).

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)

if (-not $tableShape.HasTable) {
    throw "Expected shape 2 on slide 6 to contain a table."
}

$table = $tableShape.Table
$table.ApplyStyle("{FA53ABA6-793C-4A8D-89ED-51FAC7A60CF6}")
